$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.205.15'
$ws.Range('E2').Value = '  -1.11%  '

$ws.Range('D3').Value = '2.277.92'
$ws.Range('E3').Value = '  -0.34%  '

$ws.Range('E4').Value = '  -0.47%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '112.04'
$ws.Range('E5').Value = '  +2.29%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '264.44'
$ws.Range('E6').Value = '  -0.90%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.642'
$ws.Range('E7').Value = '  +3.02%  '

$ws.Range('E8').Value = '  -0.29%  '

$ws.Range('E9').Value = '  -1.20%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '46.65'
$ws.Range('E10').Value = '  -1.35%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0937'
$ws.Range('E11').Value = '  -0.78%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '9.21'
$ws.Range('E12').Value = '  +4.88%  '

$ws.Range('E13').Value = '  +1.65%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '15.37'
$ws.Range('E14').Value = '  -1.67%  '

$ws.Range('D15').Value = '2.619.67'
$ws.Range('E15').Value = '  -0.41%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.862'
$ws.Range('E16').Value = '  +2.29%  '

$ws.Range('D17').Value = '2.268.76'
$ws.Range('E17').Value = '  -0.90%  '

$ws.Range('D18').Value = '43.192.74'
$ws.Range('E18').Value = '  -0.85%  '

$ws.Range('E19').Value = '  -1.05%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.75'
$ws.Range('E20').Value = '  +2.47%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '72.18'
$ws.Range('E21').Value = '  -0.03%  '

$ws.Range('E22').Value = '  -0.19%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '234.43'
$ws.Range('E23').Value = '  +1.20%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.87'
$ws.Range('E24').Value = '  +3.59%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.37'
$ws.Range('E25').Value = '  -2.80%  '

$ws.Range('E26').Value = '  +2.02%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.35'
$ws.Range('E27').Value = '  -1.68%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '41.50'
$ws.Range('E28').Value = '  -0.34%  '

$ws.Range('E29').Value = '  -1.59%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.25'
$ws.Range('E30').Value = '  -0.36%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '173.33'
$ws.Range('E31').Value = '  -1.48%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '21.51'
$ws.Range('E32').Value = '  +0.29%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0896'
$ws.Range('E33').Value = '  -3.01%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.67'
$ws.Range('E34').Value = '  +1.62%  '

$ws.Range('E35').Value = '  +3.51%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0380'
$ws.Range('E36').Value = '  +5.60%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.68'
$ws.Range('E37').Value = '  -0.08%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.89'
$ws.Range('E38').Value = '  +3.91%  '

$ws.Range('E39').Value = '  -2.53%  '

$ws.Range('E40').Value = '  +8.47%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '14.33'
$ws.Range('E41').Value = '  +5.52%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '75.55'
$ws.Range('E42').Value = '  +6.83%  '

$ws.Range('E43').Value = '  -2.16%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '6.11'
$ws.Range('E44').Value = '  -0.92%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.999'
$ws.Range('E45').Value = '  -0.19%  '

$ws.Range('E46').Value = '  -1.96%  '

$ws.Range('E47').Value = '  +4.32%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.55'
$ws.Range('E48').Value = '  -3.34%  '

$ws.Range('E49').Value = '  -1.31%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '100.43'
$ws.Range('E50').Value = '  -0.73%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.435'
$ws.Range('E51').Value = '  -2.29%  '
